$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 146, pushing existing rows 146:160 down to 147:161
$ws.Rows.Item(146).Insert()

# Populate the newly inserted row 146 with the new weekly record.
# Columns A,B,C,E,F,G,I,Q,R are constant across this whole data series.
$ws.Cells.Item(146, 1).Value = 9
$ws.Cells.Item(146, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(146, 3).Value = "Metropolitana"
$ws.Cells.Item(146, 4).Value = 45077
$ws.Cells.Item(146, 5).Value = 13
$ws.Cells.Item(146, 6).Value = 100112022
$ws.Cells.Item(146, 7).Value = "Arveja Verde"
$ws.Cells.Item(146, 8).Value = "Perfection"
$ws.Cells.Item(146, 9).Value = "Primera"
$ws.Cells.Item(146, 10).Value = 34
$ws.Cells.Item(146, 11).Value = 30000
$ws.Cells.Item(146, 12).Value = 32000
$ws.Cells.Item(146, 13).Value = 31000
$ws.Cells.Item(146, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(146, 15).Value = "Provincia de Huasco"
$ws.Cells.Item(146, 16).Value = 1240
$ws.Cells.Item(146, 17).Value = 25
$ws.Cells.Item(146, 18).Value = "Hortaliza"
